# Fruta / hortaliza, semanal
# Insert a brand-new weekly price record at row 8 of the data table,
# pushing the existing rows (previously 8..101) down to 9..102.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 8 (shifts rows 8:101 -> 9:102)
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new record
$ws.Range("A8").Value = 10
$ws.Range("B8").Value = 'Vega Modelo de Temuco'
$ws.Range("C8").Value = 'La Araucanía'
$ws.Range("D8").Value = 45043
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = 300000001
$ws.Range("G8").Value = 'Rabanito'
$ws.Range("H8").Value = 'Sin especificar'
$ws.Range("I8").Value = 'Primera'
$ws.Range("J8").Value = 40
$ws.Range("K8").Value = 8000
$ws.Range("L8").Value = 8000
$ws.Range("M8").Value = 8000
$ws.Range("N8").Value = '$/docena de paquetes'
$ws.Range("O8").Value = 'Provincia de Cautín'
$ws.Range("P8").Value = 667
$ws.Range("Q8").Value = 12
$ws.Range("R8").Value = 'Hortaliza'
